# Refresh crypto price/volume figures (GitHub Actions data update).
# The sheet stores Price (D) and Volume(1h) (E) as plain text, so numeric-
# looking Price values are written with a leading apostrophe to force text
# and keep the exact formatting (trailing zeros, etc.) Excel would otherwise
# normalize away when coercing a numeric-looking string to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '40.734.36'
$ws.Range("E2").Value = '  -2.29%  '

# Row 3
$ws.Range("D3").Value = '2.377.03'
$ws.Range("E3").Value = '  -3.95%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '''311.66'
$ws.Range("E5").Value = '  -2.37%  '

# Row 6
$ws.Range("D6").Value = '''86.95'
$ws.Range("E6").Value = '  -5.99%  '

# Row 7
$ws.Range("E7").Value = '  -3.68%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").Value = '''0.493'
$ws.Range("E9").Value = '  -4.18%  '

# Row 10
$ws.Range("D10").Value = '''0.0843'
$ws.Range("E10").Value = '  -3.01%  '

# Row 11
$ws.Range("D11").Value = '''30.41'
$ws.Range("E11").Value = '  -8.45%  '

# Row 12
$ws.Range("E12").Value = '  -0.50%  '

# Row 13
$ws.Range("D13").Value = '2.743.33'
$ws.Range("E13").Value = '  -3.93%  '

# Row 14
$ws.Range("E14").Value = '  -5.01%  '

# Row 15
$ws.Range("D15").Value = '''15.02'
$ws.Range("E15").Value = '  -3.44%  '

# Row 16
$ws.Range("D16").Value = '2.391.79'
$ws.Range("E16").Value = '  -3.41%  '

# Row 17
$ws.Range("E17").Value = '  -4.05%  '

# Row 18
$ws.Range("D18").Value = '40.646.68'
$ws.Range("E18").Value = '  -2.31%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0912'
$ws.Range("E19").Value = '  -3.60%  '

# Row 20
$ws.Range("D20").Value = '''6.15'
$ws.Range("E20").Value = '  -4.84%  '

# Row 21
$ws.Range("D21").Value = '''68.47'
$ws.Range("E21").Value = '  -3.20%  '

# Row 22
$ws.Range("D22").Value = '''10.78'
$ws.Range("E22").Value = '  -4.52%  '

# Row 23
$ws.Range("D23").Value = '''235.38'
$ws.Range("E23").Value = '  -1.97%  '

# Row 24
$ws.Range("E24").Value = '  -5.84%  '

# Row 25
$ws.Range("E25").Value = '  +0.07%  '

# Row 26
$ws.Range("D26").Value = '''1.80'
$ws.Range("E26").Value = '  -7.67%  '

# Row 27
$ws.Range("D27").Value = '''23.77'
$ws.Range("E27").Value = '  -4.60%  '

# Row 28
$ws.Range("E28").Value = '  -3.68%  '

# Row 29
$ws.Range("E29").Value = '  -4.97%  '

# Row 30
$ws.Range("D30").Value = '''34.30'
$ws.Range("E30").Value = '  -5.97%  '

# Row 31
$ws.Range("D31").Value = '''154.12'
$ws.Range("E31").Value = '  -1.77%  '

# Row 32
$ws.Range("E32").Value = '  -0.10%  '

# Row 33
$ws.Range("D33").Value = '''5.19'
$ws.Range("E33").Value = '  -4.95%  '

# Row 34
$ws.Range("E34").Value = '  -4.80%  '

# Row 36
$ws.Range("E36").Value = '  -2.25%  '

# Row 37
$ws.Range("E37").Value = '  -3.64%  '

# Row 38
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '''0.100'
$ws.Range("E38").Value = '  -4.12%  '

# Row 39
$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").Value = '''15.87'
$ws.Range("E39").Value = '  -8.20%  '

# Row 40
$ws.Range("E40").Value = '  -7.77%  '

# Row 41
$ws.Range("D41").Value = '''3.84'
$ws.Range("E41").Value = '  -4.58%  '

# Row 42
$ws.Range("E42").Value = '  -3.66%  '

# Row 43
$ws.Range("D43").Value = '1.965.49'
$ws.Range("E43").Value = '  -1.27%  '

# Row 44
$ws.Range("E44").Value = '  -5.23%  '

# Row 45
$ws.Range("D45").Value = '''17.80'
$ws.Range("E45").Value = '  -5.57%  '

# Row 46
$ws.Range("D46").Value = '''9.33'
$ws.Range("E46").Value = '  -1.58%  '

# Row 47
$ws.Range("E47").Value = '  -9.28%  '

# Row 48
$ws.Range("D48").Value = '2.601.37'
$ws.Range("E48").Value = '  -4.10%  '

# Row 49
$ws.Range("D49").Value = '''93.21'
$ws.Range("E49").Value = '  -5.03%  '

# Row 50
$ws.Range("D50").Value = '''71.82'
$ws.Range("E50").Value = '  -5.04%  '

# Row 51
$ws.Range("D51").Value = '''50.61'
$ws.Range("E51").Value = '  -3.40%  '
